# Daten aktualisiert am 2023-12-01
# Append three new ticker rows to the end of the existing data range.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("A183").Value = "IMX-USD"
$ws.Range("A184").Value = "TAO-USD"
$ws.Range("A185").Value = "GRT-USD"
